$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 140
$ws.Cells.Item(3, 6).Value = 391
$ws.Cells.Item(4, 6).Value = 200
$ws.Cells.Item(5, 6).Value = 33
$ws.Cells.Item(6, 6).Value = 1252
$ws.Cells.Item(7, 6).Value = 458
$ws.Cells.Item(9, 6).Value = 203
$ws.Cells.Item(10, 6).Value = 158
$ws.Cells.Item(11, 6).Value = 184
$ws.Cells.Item(12, 6).Value = 1057
$ws.Cells.Item(15, 6).Value = 196
$ws.Cells.Item(16, 6).Value = 1525
$ws.Cells.Item(17, 6).Value = 560
$ws.Cells.Item(18, 6).Value = 238
$ws.Cells.Item(19, 6).Value = 356
$ws.Cells.Item(21, 6).Value = 846
$ws.Cells.Item(22, 6).Value = 1162
$ws.Cells.Item(25, 6).Value = 2678
$ws.Cells.Item(26, 6).Value = 1467
$ws.Cells.Item(27, 6).Value = 68
$ws.Cells.Item(28, 6).Value = 47
$ws.Cells.Item(29, 6).Value = 441
$ws.Cells.Item(30, 6).Value = 514
$ws.Cells.Item(31, 6).Value = 1304
$ws.Cells.Item(32, 6).Value = 834
$ws.Cells.Item(33, 6).Value = 1403
$ws.Cells.Item(34, 6).Value = 166
$ws.Cells.Item(36, 6).Value = 794
$ws.Cells.Item(37, 6).Value = 642
$ws.Cells.Item(38, 6).Value = 690
$ws.Cells.Item(39, 6).Value = 871
$ws.Cells.Item(40, 6).Value = 373
$ws.Cells.Item(41, 6).Value = 259
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 180
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(15, 6).Value = 645
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(6, 6).Value = 140
$ws.Cells.Item(7, 6).Value = 391
$ws.Cells.Item(8, 6).Value = 200
$ws.Cells.Item(9, 6).Value = 33
$ws.Cells.Item(11, 6).Value = 180
$ws.Cells.Item(12, 6).Value = 1252
$ws.Cells.Item(13, 6).Value = 458
$ws.Cells.Item(15, 6).Value = 203
$ws.Cells.Item(17, 6).Value = 158
$ws.Cells.Item(18, 6).Value = 184
$ws.Cells.Item(19, 6).Value = 1057
$ws.Cells.Item(21, 6).Value = 196
$ws.Cells.Item(22, 6).Value = 1525
$ws.Cells.Item(23, 6).Value = 560
$ws.Cells.Item(24, 6).Value = 238
$ws.Cells.Item(25, 6).Value = 356
$ws.Cells.Item(26, 6).Value = 3
$ws.Cells.Item(29, 6).Value = 1162
$ws.Cells.Item(30, 6).Value = 2678
$ws.Cells.Item(31, 6).Value = 1467
$ws.Cells.Item(32, 6).Value = 68
$ws.Cells.Item(35, 6).Value = 441
$ws.Cells.Item(36, 6).Value = 1304
$ws.Cells.Item(39, 6).Value = 834
$ws.Cells.Item(40, 6).Value = 1403
$ws.Cells.Item(41, 6).Value = 794
$ws.Cells.Item(42, 6).Value = 642
$ws.Cells.Item(43, 6).Value = 690
$ws.Cells.Item(44, 6).Value = 871
$ws.Cells.Item(45, 6).Value = 373
$ws.Cells.Item(48, 6).Value = 259
